# Insert a new data row at row 195 (pushes existing rows 195-212 down to 196-213)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("195:195").Insert()

# Fill in the new record on row 195
$ws.Cells.Item(195, 1).Value = 8
$ws.Cells.Item(195, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(195, 3).Value = "Coquimbo"
$ws.Cells.Item(195, 4).Value = 44578
$ws.Cells.Item(195, 5).Value = 4
$ws.Cells.Item(195, 6).Value = 100112012
$ws.Cells.Item(195, 7).Value = "Espinaca"
$ws.Cells.Item(195, 8).Value = "Sin especificar"
$ws.Cells.Item(195, 9).Value = "Primera"
$ws.Cells.Item(195, 10).Value = 2400
$ws.Cells.Item(195, 11).Value = 400
$ws.Cells.Item(195, 12).Value = 500
$ws.Cells.Item(195, 13).Value = 450
$ws.Cells.Item(195, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(195, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(195, 16).Value = 900
$ws.Cells.Item(195, 17).Value = 0.5
$ws.Cells.Item(195, 18).Value = "Hortaliza"

# Keep the numeric/date format (s="2") used by column D on the other rows
$ws.Cells.Item(195, 4).NumberFormat = $ws.Cells.Item(196, 4).NumberFormat
